# fix(import): update date-time format & edit data
#
# template_stok.xlsx: the "stok_tanggal" (stock date) column is dropped
# from the template - header text in D1 and its sample value in D2 are
# both removed, and the active selection moves from H6 to F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "stok_tanggal" header text from D1 but keep the cell's
# (bold-header) formatting in place.
$ws.Range("D1").ClearContents()

# Remove the sample date value from D2 entirely (value + number format),
# leaving the cell completely blank.
$ws.Range("D2").Clear()

# Move the selected/active cell to F4.
$ws.Range("F4").Select()
